$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Metadata sheet: bump the "Date" property to the new timestamp
# ------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-19T13:17:15+00:00"

# ------------------------------------------------------------------
# 2) Elements sheet: the two mapping columns (AK = "RIM Mapping",
#    AL = "Spécification métier...") were swapped - the business
#    mapping column now comes first, RIM Mapping second. Swap both
#    the cell contents (header + data rows) and the column widths.
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elements")

$lastRow = 6
for ($r = 1; $r -le $lastRow; $r++) {
    $akCell = $ws.Cells.Item($r, 37)
    $alCell = $ws.Cells.Item($r, 38)

    $akVal = $akCell.Value2
    $alVal = $alCell.Value2

    if ($akVal -eq $null) { $akVal = "" }
    if ($alVal -eq $null) { $alVal = "" }

    # Skip cells that already hold identical content on both sides -
    # nothing changes, so leave the original cell (and its shared
    # string typing) completely untouched.
    if ($akVal -ceq $alVal) { continue }

    $akCell.Value = $alVal
    $alCell.Value = $akVal
}

# Swap the stored column widths too (AK was the narrow column, AL the
# wide one - now it's the other way round).
$ws.Range("AK1").ColumnWidth = 59.0
$ws.Range("AL1").ColumnWidth = 24.15
